$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 520.4
$ws.Range("I2").Value = 700.3333
$ws.Range("K2").Value = 700.3333
$ws.Range("M2").Value = -587.3333
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -230
$ws.Range("H53").Value = 745.7143
$ws.Range("I53").Value = 582.75
$ws.Range("J53").Value = 810.9
$ws.Range("K53").Value = 582.75
$ws.Range("L53").Value = 810.9
$ws.Range("M53").Value = 54.25
$ws.Range("N53").Value = -2084.9
$ws.Range("H88").Value = 5220
$ws.Range("J88").Value = 3931.2222
$ws.Range("L88").Value = 3931.2222
$ws.Range("N88").Value = -4743.2222
$ws.Range("H91").Value = 5220
$ws.Range("J91").Value = 3931.2222
$ws.Range("L91").Value = 3931.2222
$ws.Range("N91").Value = -6739.2222
$ws.Range("H112").Value = 1309.3214
$ws.Range("I112").Value = 1333.3334
$ws.Range("J112").Value = 1306.44
$ws.Range("K112").Value = 4000.0002
$ws.Range("L112").Value = 3919.32
$ws.Range("M112").Value = -2892.0002
$ws.Range("N112").Value = -6135.32
$ws.Range("H132").Value = 25613.363
$ws.Range("I132").Value = 1906.125
$ws.Range("K132").Value = 5718.375
$ws.Range("M132").Value = -3188.375
$ws.Range("H135").Value = 3192.2334
$ws.Range("J135").Value = 4326.909
$ws.Range("L135").Value = 38942.181
$ws.Range("N135").Value = -44012.181
$ws.Range("H138").Value = 3052.85
$ws.Range("I138").Value = 1436.75
$ws.Range("K138").Value = 4310.25
$ws.Range("M138").Value = 829.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 37800
$ws.Range("I28").Value = 14400
$ws.Range("J28").Value = 49500
$ws.Range("K28").Value = 14400
$ws.Range("L28").Value = 49500
$ws.Range("M28").Value = -14208
$ws.Range("N28").Value = -49884
$ws.Range("H61").Value = 9668.5
$ws.Range("I61").Value = 1464.4546
$ws.Range("K61").Value = 1464.4546
$ws.Range("M61").Value = -1252.4546
$ws.Range("H99").Value = 37800
$ws.Range("I99").Value = 14400
$ws.Range("J99").Value = 49500
$ws.Range("K99").Value = 14400
$ws.Range("L99").Value = 49500
$ws.Range("M99").Value = -11405
$ws.Range("N99").Value = -55490
$ws.Range("H136").Value = 9668.5
$ws.Range("I136").Value = 1464.4546
$ws.Range("K136").Value = 4393.3638
$ws.Range("M136").Value = -1843.3638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5389.4736
$ws.Range("I20").Value = 3809.182
$ws.Range("K20").Value = 3809.182
$ws.Range("M20").Value = -3562.182
$ws.Range("H87").Value = 48809.477
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 48809.477
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H94").Value = 2927758
$ws.Range("I94").Value = 3321.037
$ws.Range("K94").Value = 3321.037
$ws.Range("M94").Value = -2870.037
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2007.7307
$ws.Range("I132").Value = 1763.0454
$ws.Range("J132").Value = 3353.5
$ws.Range("K132").Value = 5289.1362
$ws.Range("L132").Value = 10060.5
$ws.Range("M132").Value = -2759.1362
$ws.Range("N132").Value = -15120.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 913.3333
$ws.Range("I23").Value = 2223
$ws.Range("J23").Value = 409.6154
$ws.Range("K23").Value = 6669
$ws.Range("L23").Value = 1228.8462
$ws.Range("M23").Value = -6434
$ws.Range("N23").Value = -1698.8462
$ws.Range("H26").Value = 429.83334
$ws.Range("I26").Value = 15.8
$ws.Range("K26").Value = 47.40000000000001
$ws.Range("M26").Value = 240.6
$ws.Range("H68").Value = 2065.7778
$ws.Range("I68").Value = 844.5
$ws.Range("K68").Value = 2533.5
$ws.Range("M68").Value = -1722.5
$ws.Range("H71").Value = 2065.7778
$ws.Range("I71").Value = 844.5
$ws.Range("K71").Value = 7600.5
$ws.Range("M71").Value = -3544.5
$ws.Range("H113").Value = 1123.8462
$ws.Range("I113").Value = 2370
$ws.Range("J113").Value = 897.2727
$ws.Range("K113").Value = 7110
$ws.Range("L113").Value = 2691.8181
$ws.Range("M113").Value = -4940
$ws.Range("N113").Value = -7031.8181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 76840.664
$ws.Range("J39").Value = 76840.664
$ws.Range("L39").Value = 76840.664
$ws.Range("N39").Value = -77904.664
$ws.Range("H49").Value = 23045
$ws.Range("J49").Value = 36495
$ws.Range("L49").Value = 36495
$ws.Range("N49").Value = -36863
$ws.Range("H80").Value = 25644714
$ws.Range("I80").Value = 1657.8334
$ws.Range("K80").Value = 1657.8334
$ws.Range("M80").Value = -659.8334
$ws.Range("H83").Value = 25644714
$ws.Range("I83").Value = 1657.8334
$ws.Range("K83").Value = 8289.166999999999
$ws.Range("M83").Value = -3297.166999999999
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H96").Value = 17630.5
$ws.Range("J96").Value = 17630.5
$ws.Range("L96").Value = 17630.5
$ws.Range("N96").Value = -23122.5
$ws.Range("H98").Value = 5000
$ws.Range("J98").Value = 5000
$ws.Range("L98").Value = 5000
$ws.Range("N98").Value = -10990
$ws.Range("H99").Value = 28470.5
$ws.Range("I99").Value = 16681.5
$ws.Range("K99").Value = 16681.5
$ws.Range("M99").Value = -14435.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 31251514
$ws.Range("I16").Value = 83334280
$ws.Range("K16").Value = 83334280
$ws.Range("M16").Value = -83334110
$ws.Range("H43").Value = 35795
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 3631.5557
$ws.Range("J46").Value = 3898
$ws.Range("L46").Value = 3898
$ws.Range("N46").Value = -4274
$ws.Range("H61").Value = 1845.6875
$ws.Range("I61").Value = 1836.1428
$ws.Range("K61").Value = 1836.1428
$ws.Range("M61").Value = -1634.1428
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H106").Value = 11705
$ws.Range("J106").Value = 11705
$ws.Range("L106").Value = 11705
$ws.Range("N106").Value = -14229
$ws.Range("H113").Value = 1845.6875
$ws.Range("I113").Value = 1836.1428
$ws.Range("K113").Value = 1836.1428
$ws.Range("M113").Value = 333.8571999999999
$ws.Range("H122").Value = 3661
$ws.Range("I122").Value = 3661
$ws.Range("K122").Value = 10983
$ws.Range("M122").Value = -8533
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 49950
$ws.Range("J101").Value = 49950
$ws.Range("L101").Value = 49950
$ws.Range("N101").Value = -56440
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H141").Value = 82391
$ws.Range("J141").Value = 82391
$ws.Range("L141").Value = 82391
$ws.Range("N141").Value = -92751
